$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column I (Segment) ---
# Target OOXML width is 50.68359375 characters. The COM ColumnWidth setter in
# this host rounds to whole pixels with a +10px pad (px = round(chars*12+10)),
# so we pre-compensate the input to land as close as possible on the target.
$ws.Columns("I").ColumnWidth = 49.833333333333336

# --- Append two new coded-segment rows (142 and 143) ---
# Clone formatting (fill/border/font/number-format) from the last existing
# data row so the new rows render identically to the rest of the table.
$ws.Range("A141:M141").Copy()
$ws.Range("A142:M143").PasteSpecial(-4122)

# Row 142
$ws.Range("A142").Value = "\u25cf"
$ws.Range("D142").Value = "10323"
$ws.Range("E142").Value = "Location:Country"
$ws.Range("F142").Value = "2: 3521"
$ws.Range("G142").Value = "2: 3525"
$ws.Range("H142").Value = 0
$ws.Range("I142").Value = "Japan"
$ws.Range("J142").Value = 5
$ws.Range("K142").Value = 0.037268932617769826
$ws.Range("L142").Value = "dattaray"
$ws.Range("M142").Value = "9/17/2019 12:14:42"
$ws.Rows("142").RowHeight = 15.6

# Row 143
$ws.Range("A143").Value = "\u25cf"
$ws.Range("D143").Value = "10323"
$ws.Range("E143").Value = "Exclusion: Not first case in country"
$ws.Range("F143").Value = "3: 424"
$ws.Range("G143").Value = "3: 497"
$ws.Range("H143").Value = 0
$ws.Range("I143").Value = "To our knowledge, this is `nthe first report of a fatality caused by ISMRK."
$ws.Range("J143").Value = 74
$ws.Range("K143").Value = 0.55158020274299346
$ws.Range("L143").Value = "dattaray"
$ws.Range("M143").Value = "9/17/2019 12:15:52"
$ws.Rows("143").RowHeight = 25.8
